# [Fonds de solidarite] Add 2022-06-14 data
# Update nombre_aides (column C) and montant_total (column E) for the rows
# whose underlying cumulative counters changed with the new data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;   C = 100827; E = 327352874 },
    @{ Row = 9;   C = 285;    E = 36555237 },
    @{ Row = 19;  C = 4367;   E = 66092282 },
    @{ Row = 48;  C = 1677;   E = 31738250 },
    @{ Row = 56;  C = 11976;  E = 187859672 },
    @{ Row = 64;  C = 5212;   E = 20422343 },
    @{ Row = 92;  C = 409189; E = 1595779968 },
    @{ Row = 93;  C = 209615; E = 1309459328 },
    @{ Row = 94;  C = 94218;  E = 918462438 },
    @{ Row = 95;  C = 50782;  E = 933328284 },
    @{ Row = 97;  C = 2162;   E = 214351518 },
    @{ Row = 104; C = 135253; E = 272253470 },
    @{ Row = 119; C = 356;    E = 10745878 },
    @{ Row = 173; C = 96860;  E = 327935644 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
